$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 835
$ws.Range("I4").Value = 520
$ws.Range("J4").Value = 1060
$ws.Range("K4").Value = 520
$ws.Range("L4").Value = 1060
$ws.Range("M4").Value = -406
$ws.Range("N4").Value = -1288
$ws.Range("H8").Value = 30.333334
$ws.Range("I8").Value = 25.5
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 76.5
$ws.Range("L8").Value = 120
$ws.Range("M8").Value = 62.5
$ws.Range("N8").Value = -398
$ws.Range("H15").Value = 42.5
$ws.Range("I15").Value = 42.5
$ws.Range("K15").Value = 127.5
$ws.Range("M15").Value = 41.5
$ws.Range("H112").Value = 1115
$ws.Range("J112").Value = 1164.2858
$ws.Range("L112").Value = 3492.8574
$ws.Range("N112").Value = -5708.857400000001
$ws.Range("H125").Value = 1346.4
$ws.Range("I125").Value = 910.6667
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 8196.0003
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -5736.0003
$ws.Range("N125").Value = -22920
$ws.Range("H126").Value = 47763
$ws.Range("J126").Value = 47763
$ws.Range("L126").Value = 47763
$ws.Range("N126").Value = -57643
$ws.Range("H127").Value = 2022.2667
$ws.Range("J127").Value = 1349
$ws.Range("L127").Value = 4047
$ws.Range("N127").Value = -13967
$ws.Range("H129").Value = 401628.9
$ws.Range("J129").Value = 1504.4
$ws.Range("L129").Value = 4513.200000000001
$ws.Range("N129").Value = -14513.2
$ws.Range("H132").Value = 28893.03
$ws.Range("I132").Value = 4434.8066
$ws.Range("J132").Value = 218444.25
$ws.Range("K132").Value = 13304.4198
$ws.Range("L132").Value = 655332.75
$ws.Range("M132").Value = -10774.4198
$ws.Range("N132").Value = -660392.75
$ws.Range("H133").Value = 77339.14
$ws.Range("J133").Value = 77339.14
$ws.Range("L133").Value = 77339.14
$ws.Range("N133").Value = -87459.14

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16763.535
$ws.Range("I32").Value = 16084
$ws.Range("J32").Value = 19899.846
$ws.Range("K32").Value = 16084
$ws.Range("L32").Value = 19899.846
$ws.Range("M32").Value = -15797
$ws.Range("N32").Value = -20473.846
$ws.Range("H109").Value = 43368.668
$ws.Range("J109").Value = 43368.668
$ws.Range("L109").Value = 43368.668
$ws.Range("N109").Value = -46142.668
$ws.Range("H112").Value = 7165396.5
$ws.Range("J112").Value = 7165396.5
$ws.Range("L112").Value = 7165396.5
$ws.Range("N112").Value = -7168350.5
$ws.Range("H113").Value = 41247.285
$ws.Range("J113").Value = 41247.285
$ws.Range("L113").Value = 41247.285
$ws.Range("N113").Value = -49925.285

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 36095
$ws.Range("J35").Value = 36095
$ws.Range("L35").Value = 36095
$ws.Range("N35").Value = -36715
$ws.Range("H42").Value = 180000
$ws.Range("J42").Value = 180000
$ws.Range("L42").Value = 180000
$ws.Range("N42").Value = -180656

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 52606.5
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 52606.5
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 52606.5
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -53196.5
$ws.Range("H34").Value = 52606.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 52606.5
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 52606.5
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -53010.5
$ws.Range("H132").Value = 32653.088
$ws.Range("I132").Value = 1679.3636
$ws.Range("J132").Value = 111278.695
$ws.Range("K132").Value = 5038.0908
$ws.Range("L132").Value = 333836.085
$ws.Range("M132").Value = -2508.0908
$ws.Range("N132").Value = -338896.085

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H129").Value = 105007.27
$ws.Range("I129").Value = 231359.23
$ws.Range("J129").Value = 2346.3125
$ws.Range("K129").Value = 694077.6900000001
$ws.Range("L129").Value = 7038.9375
$ws.Range("M129").Value = -689077.6900000001
$ws.Range("N129").Value = -17038.9375
$ws.Range("H131").Value = 42456.844
$ws.Range("J131").Value = 45814.89
$ws.Range("L131").Value = 137444.67
$ws.Range("N131").Value = -147524.67

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 43485.2
$ws.Range("J114").Value = 43485.2
$ws.Range("L114").Value = 43485.2
$ws.Range("N114").Value = -52163.2
$ws.Range("H126").Value = 7183.05
$ws.Range("I126").Value = 14601.375
$ws.Range("J126").Value = 2237.5
$ws.Range("K126").Value = 43804.125
$ws.Range("L126").Value = 6712.5
$ws.Range("M126").Value = -41334.125
$ws.Range("N126").Value = -11652.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2062
$ws.Range("I68").Value = 2100.5
$ws.Range("J68").Value = 2034
$ws.Range("K68").Value = 2100.5
$ws.Range("L68").Value = 2034
$ws.Range("M68").Value = -1351.5
$ws.Range("N68").Value = -3532
$ws.Range("H71").Value = 2062
$ws.Range("I71").Value = 2100.5
$ws.Range("J71").Value = 2034
$ws.Range("K71").Value = 10502.5
$ws.Range("L71").Value = 10170
$ws.Range("M71").Value = -6758.5
$ws.Range("N71").Value = -17658
$ws.Range("H110").Value = 30749.75
$ws.Range("J110").Value = 30749.75
$ws.Range("L110").Value = 30749.75
$ws.Range("N110").Value = -38929.75
$ws.Range("H111").Value = 38585.75
$ws.Range("J111").Value = 38585.75
$ws.Range("L111").Value = 38585.75
$ws.Range("N111").Value = -46765.75
$ws.Range("H132").Value = 3149.6956
$ws.Range("I132").Value = 2772.425
$ws.Range("J132").Value = 5664.8335
$ws.Range("K132").Value = 8317.275000000001
$ws.Range("L132").Value = 16994.5005
$ws.Range("M132").Value = -5787.275000000001
$ws.Range("N132").Value = -22054.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 35394.57
$ws.Range("J114").Value = 35394.57
$ws.Range("L114").Value = 35394.57
$ws.Range("N114").Value = -44072.57
$ws.Range("H119").Value = 48694
$ws.Range("J119").Value = 48694
$ws.Range("L119").Value = 48694
$ws.Range("N119").Value = -58370
$ws.Range("H126").Value = 4202494.5
$ws.Range("I126").Value = 4902577
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 14707731
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -14705261
$ws.Range("N126").Value = -10940

Write-Output "edits applied"